$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header values in row 1 (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 (CON) values B2:E2
$ws.Range("B2").Value = 5.1090777001051748
$ws.Range("C2").Value = 5.7056655588445722
$ws.Range("D2").Value = 3.2166496700074267
$ws.Range("E2").Value = 3.6702370794705077

# Update row 3 (STR) values B3:E3
$ws.Range("B3").Value = 6.5387440478270467
$ws.Range("C3").Value = 8.6722543582275193
$ws.Range("D3").Value = 4.4255690177005844
$ws.Range("E3").Value = 1.5497740631580796

# Update the selection on the sheet to match B1:E3
$ws.Range("B1:E3").Select()
